$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '42.809.72'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '2.336.17'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.93%  '
$ws.Range('E7').Value = '  -5.28%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.511'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.13'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0800'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.47%  '
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.82'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.82'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.91%  '
$ws.Range('D16').Value = '2.313.12'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.800'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('D18').Value = '42.728.19'
$ws.Range('E18').Value = '  -1.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('D20').Value = '0.0₃0910'
$ws.Range('E20').Value = '  -2.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.66'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.89%  '
$ws.Range('E28').Value = '  +6.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.91'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.67'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.98%  '
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.12'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.88%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.60'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.92%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.39'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0728'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.57%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.45'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.48%  '
$ws.Range('E38').Value = '  -5.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.87'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('E40').Value = '  -3.73%  '
$ws.Range('E41').Value = '  -3.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.35'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Value = '2.029.39'
$ws.Range('E43').Value = '  +2.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0286'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.61'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.93'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '56.36'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.94%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('D50').Value = '2.560.49'
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('E51').Value = '  +1.47%  '
